# Insert a new observation row ("KB1-2") before the current row 68,
# pushing all subsequent rows (old 68..145) down by one (new 69..146).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(68).Insert() | Out-Null

$ws.Range("A68").Value = "KB1-2"
$ws.Range("B68").Value = 427404.96500000003
$ws.Range("C68").Value = 6654944.7280000001
$ws.Range("D68").Value = 2

# Update the view so the same cell is selected / visible as in the
# authored workbook.
$ws.Range("F69").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
